# Append a new trade row (row 6) to the sheet, matching the columns:
# A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice, E=IsShortSell,
# F=Price Change %, G=Date, H=Profitable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 9950.73
$ws.Range("B6").Value = 10027.950000000001
$ws.Range("C6").Value = 79.650000000000006
$ws.Range("D6").Value = 79.040000000000006
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -0.77
$ws.Range("G6").Value = 42612.67423611111
$ws.Range("H6").Value = $false

# Give the new Date cell (G6) the same formatting (date/time number format)
# as the cell above it, without introducing a brand-new style entry.
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
